$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.055778869259484
$ws.Range("C2").Value = 0.2570746161543127
$ws.Range("E2").Value = 0.1106345412705487
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002431876339472617
$ws.Range("I2").Value = 0.6340566669494354
$ws.Range("L2").Value = 0.2027027211253767
$ws.Range("M2").Value = 0.2270515722764372
$ws.Range("O2").Value = 2.340885380319065
$ws.Range("B3").Value = 0.9464823195150984
$ws.Range("C3").Value = 0.2409517690229563
$ws.Range("E3").Value = 0.1119467971376099
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002434498755688357
$ws.Range("I3").Value = 0.6467971415717475
$ws.Range("L3").Value = 0.2002138771724162
$ws.Range("M3").Value = 0.2101091209853223
$ws.Range("O3").Value = 2.376906891880083
$ws.Range("B4").Value = 0.8792733638976529
$ws.Range("C4").Value = 0.2310049221545114
$ws.Range("E4").Value = 0.1127997913745087
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002436194233340239
$ws.Range("I4").Value = 0.6551167120367491
$ws.Range("L4").Value = 0.1987883573293345
$ws.Range("M4").Value = 0.1997324592189145
$ws.Range("O4").Value = 2.401040341550555
$ws.Range("B5").Value = 0.8518617275983047
$ws.Range("C5").Value = 0.2269398774795377
$ws.Range("E5").Value = 0.1131592932593644
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002436906667302348
$ws.Range("I5").Value = 0.6586317769289725
$ws.Range("L5").Value = 0.1982333104203491
$ws.Range("M5").Value = 0.1955107156197187
$ws.Range("O5").Value = 2.411380984276036
$ws.Range("B6").Value = 0.8473086841531767
$ws.Range("C6").Value = 0.2262641862538999
$ws.Range("E6").Value = 0.1132197074314957
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002437026267640393
$ws.Range("I6").Value = 0.6592229817616673
$ws.Range("L6").Value = 0.1981427090285237
$ws.Range("M6").Value = 0.1948101189639218
$ws.Range("O6").Value = 2.413128578903724
$ws.Range("B7").Value = 0.8789037734827616
$ws.Range("C7").Value = 0.2309501462268315
$ws.Range("E7").Value = 0.1128045915408518
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002436203754209035
$ws.Range("I7").Value = 0.6551636124576241
$ws.Range("L7").Value = 0.1987807669776842
$ws.Range("M7").Value = 0.1996754952873943
$ws.Range("O7").Value = 2.401177751227934
$ws.Range("B8").Value = 1.018115313472435
$ws.Range("C8").Value = 0.2515254512938725
$ws.Range("E8").Value = 0.1110772065700472
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002432762882256531
$ws.Range("I8").Value = 0.6383463576580954
$ws.Range("L8").Value = 0.2018232958969932
$ws.Range("M8").Value = 0.2212045624477952
$ws.Range("O8").Value = 2.352886619685847
$ws.Range("B9").Value = 1.290245889010123
$ws.Range("C9").Value = 0.2914873296000735
$ws.Range("E9").Value = 0.1080641227255026
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002426689286847764
$ws.Range("I9").Value = 0.6093173066333346
$ws.Range("L9").Value = 0.2086022661259719
$ws.Range("M9").Value = 0.2636200276306795
$ws.Range("O9").Value = 2.27422473689623
$ws.Range("B10").Value = 1.489585780646962
$ws.Range("C10").Value = 0.3206006270932278
$ws.Range("E10").Value = 0.106077582624354
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002422633799250487
$ws.Range("I10").Value = 0.5904065188292584
$ws.Range("L10").Value = 0.214076516681601
$ws.Range("M10").Value = 0.2948928246197582
$ws.Range("O10").Value = 2.226263203806852
$ws.Range("B11").Value = 1.580128523921758
$ws.Range("C11").Value = 0.3337892827826181
$ws.Range("E11").Value = 0.1052229644850244
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.00242087631482086
$ws.Range("I11").Value = 0.582330238005996
$ws.Range("L11").Value = 0.2166738162984956
$ws.Range("M11").Value = 0.3091416668067808
$ws.Range("O11").Value = 2.206590795029243
$ws.Range("B12").Value = 1.614393348922306
$ws.Range("C12").Value = 0.3387753150733772
$ws.Range("E12").Value = 0.1049063834889141
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002420223300755642
$ws.Range("I12").Value = 0.5793478264756171
$ws.Range("L12").Value = 0.2176726964253533
$ws.Range("M12").Value = 0.3145403689050568
$ws.Range("O12").Value = 2.199450873375156
$ws.Range("B13").Value = 1.607014803163224
$ws.Range("C13").Value = 0.3377018543489214
$ws.Range("E13").Value = 0.1049742517571805
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002420363383557649
$ws.Range("I13").Value = 0.5799867630353184
$ws.Range("L13").Value = 0.2174568883398109
$ws.Range("M13").Value = 0.313377534070348
$ws.Range("O13").Value = 2.200974797082438
$ws.Range("B14").Value = 1.582947958398847
$ws.Range("C14").Value = 0.3341996532068094
$ws.Range("E14").Value = 0.1051967780634691
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002420822340689031
$ws.Range("I14").Value = 0.5820833508204437
$ws.Range("L14").Value = 0.216755687648714
$ws.Range("M14").Value = 0.3095857633000563
$ws.Range("O14").Value = 2.205997178580276
$ws.Range("B15").Value = 1.568203432678899
$ws.Range("C15").Value = 0.3320533739011751
$ws.Range("E15").Value = 0.1053339988872288
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002421105092188347
$ws.Range("I15").Value = 0.5833774622678138
$ws.Range("L15").Value = 0.2163281780316169
$ws.Range("M15").Value = 0.3072635731148949
$ws.Range("O15").Value = 2.209113881975938
$ws.Range("B16").Value = 1.483665659722192
$ws.Range("C16").Value = 0.3197375825571385
$ws.Range("E16").Value = 0.1061344203192287
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002422750408216102
$ws.Range("I16").Value = 0.5909449322859146
$ws.Range("L16").Value = 0.2139089263856278
$ws.Range("M16").Value = 0.2939620611852973
$ws.Range("O16").Value = 2.227592106947966
$ws.Range("B17").Value = 1.431767774301306
$ws.Range("C17").Value = 0.3121679107037494
$ws.Range("E17").Value = 0.1066380123136043
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002423782094551455
$ws.Range("I17").Value = 0.5957222891298564
$ws.Range("L17").Value = 0.2124521708753946
$ws.Range("M17").Value = 0.285807620301469
$ws.Range("O17").Value = 2.239478182848174
$ws.Range("B18").Value = 1.401904609482756
$ws.Range("C18").Value = 0.3078088642560317
$ws.Range("E18").Value = 0.1069322844025064
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002424383721010998
$ws.Range("I18").Value = 0.5985196263957846
$ws.Range("L18").Value = 0.2116243614589308
$ws.Range("M18").Value = 0.281119554347022
$ws.Range("O18").Value = 2.246516623495324
$ws.Range("B19").Value = 1.391791302269723
$ws.Range("C19").Value = 0.3063320869537733
$ws.Range("E19").Value = 0.1070327136017124
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002424588836296978
$ws.Range("I19").Value = 0.5994752566543333
$ws.Range("L19").Value = 0.2113458118907516
$ws.Range("M19").Value = 0.2795326343171993
$ws.Range("O19").Value = 2.248934359774822
$ws.Range("B20").Value = 1.437293739049778
$ws.Range("C20").Value = 0.3129742525840982
$ws.Range("E20").Value = 0.1065839260916675
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002423671418757985
$ws.Range("I20").Value = 0.5952086038208435
$ws.Range("L20").Value = 0.2126062023501163
$ws.Range("M20").Value = 0.2866754537443086
$ws.Range("O20").Value = 2.238191987773789
$ws.Range("B21").Value = 1.590017579411722
$ws.Range("C21").Value = 0.3352285601714584
$ws.Range("E21").Value = 0.1051312256045356
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002420687194822532
$ws.Range("I21").Value = 0.5814654705298885
$ws.Range("L21").Value = 0.216961231406728
$ws.Range("M21").Value = 0.3106994198382651
$ws.Range("O21").Value = 2.204513572647812
$ws.Range("B22").Value = 1.689703807165699
$ws.Range("C22").Value = 0.3497249565424454
$ws.Range("E22").Value = 0.1042228563243973
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002418809706001034
$ws.Range("I22").Value = 0.5729260456298935
$ws.Range("L22").Value = 0.2198968785959181
$ws.Range("M22").Value = 0.3264176777614551
$ws.Range("O22").Value = 2.184307816676366
$ws.Range("B23").Value = 1.636511695268553
$ws.Range("C23").Value = 0.3419924540014847
$ws.Range("E23").Value = 0.1047039174635833
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002419805108159869
$ws.Range("I23").Value = 0.5774431370242539
$ws.Range("L23").Value = 0.2183219071018243
$ws.Range("M23").Value = 0.3180270678168213
$ws.Range("O23").Value = 2.194926496387623
$ws.Range("B24").Value = 1.434795532412068
$ws.Range("C24").Value = 0.3126097276441726
$ws.Range("E24").Value = 0.1066083636749189
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002423721428675867
$ws.Range("I24").Value = 0.5954406828150844
$ws.Range("L24").Value = 0.2125365345063699
$ws.Range("M24").Value = 0.2862831061276765
$ws.Range("O24").Value = 2.23877283815284
$ws.Range("B25").Value = 1.216726868874048
$ws.Range("C25").Value = 0.280719104100541
$ws.Range("E25").Value = 0.1088392745947201
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002428260627489754
$ws.Range("I25").Value = 0.6167466667740822
$ws.Range("L25").Value = 0.2066815518849836
$ws.Range("M25").Value = 0.2521254653495433
$ws.Range("O25").Value = 2.293782320316595
